# feat!: removal of option `fieldMatchType`
# new default is `labelTypeBrackets` to avoid any collisions with labels and types
#
# Rewrite the header row (A1:H1) so every column label is suffixed with its
# technical field name in brackets, e.g. "ID" -> "ID[product_ID]".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A1").Value = "ID[product_ID]"
$ws.Range("B1").Value = "Quantity[quantity]"
$ws.Range("C1").Value = "ProductTitle[title]"
$ws.Range("D1").Value = "UnitPrice[price]"
$ws.Range("E1").Value = "validFrom[validFrom]"
$ws.Range("F1").Value = "timestamp[timestamp]"
$ws.Range("G1").Value = "date[date]"
$ws.Range("H1").Value = "time[time]"

# Match the header-row selection captured in the fixture (A1:H1).
$ws.Range("A1:H1").Select()
